# Update the roster table on Sheet1 (A2:C19) to reflect the new
# player / position / team assignments.
#
# Column A = Oyuncu Adı (Player name)
# Column B = Pozisyon (Position)
# Column C = Takım (Team)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row=2;  A="Ja Morant";        B="PG";          C="Memphis Grizzlies" },
    @{ Row=3;  A="Isaiah Collier";   B="PG,SG";       C="Utah Jazz" },
    @{ Row=4;  A="Max Strus";        B="SG,SF";       C="Cleveland Cavaliers" },
    @{ Row=5;  A="Tyler Herro";      B="PG,SG";       C="Miami Heat" },
    @{ Row=6;  A="DeMar DeRozan";    B="SF,PF";       C="Sacramento Kings" },
    @{ Row=7;  A="Miles Bridges";    B="SF,PF";       C="Charlotte Hornets" },
    @{ Row=8;  A="Christian Braun";  B="SG,SF";       C="Denver Nuggets" },
    @{ Row=9;  A="Evan Mobley";      B="PF,C";        C="Cleveland Cavaliers" },
    @{ Row=10; A="Brook Lopez";      B="C";           C="Milwaukee Bucks" },
    @{ Row=11; A="Max Christie";     B="SG,SF";       C="Dallas Mavericks" },
    @{ Row=12; A="Luka Doncic";      B="PG,SG";       C="Los Angeles Lakers" },
    @{ Row=13; A="Zach Collins";     B="PF,C";        C="Chicago Bulls" },
    @{ Row=14; A="Scottie Barnes";   B="PG,SG,SF,PF"; C="Toronto Raptors" },
    @{ Row=15; A="De'Aaron Fox";     B="PG,SG";       C="San Antonio Spurs" },
    @{ Row=16; A="Mikal Bridges";    B="SG,SF,PF";    C="New York Knicks" },
    @{ Row=17; A="Nikola Vucevic";   B="PF,C";        C="Chicago Bulls" },
    @{ Row=18; A="P.J. Washington";  B="SF,PF";       C="Dallas Mavericks" },
    @{ Row=19; A="Josh Giddey";      B="PG,SG,SF";    C="Chicago Bulls" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}
